# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the second handback file (c0d0d805-...) on the per-language
# sheets, and the rolled-up "Latest HO Xliff Generate Date" on the Overview
# sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the c0d0d805-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-20 00:23:28"

# --- zh-cn sheet: row 3 is the c0d0d805-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-20 00:23:17"
$wsZhCn.Range("K3").Value = "2016-10-20 00:24:00"

# --- de-de sheet: row 3 is the c0d0d805-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-10-20 00:23:28"
$wsDeDe.Range("K3").Value = "2016-10-20 00:24:20"
